$d = $word.ActiveDocument
$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs($count)
$r = $lastPara.Range
$xml = '<w:p><w:r><w:rPr/><w:t>Having never worked with HTML and CSS</w:t></w:r><w:r><w:rPr/><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr/><w:t>I</w:t></w:r><w:r><w:rPr/><w:t xml:space="preserve"> struggled at the very </w:t></w:r><w:r><w:rPr/><w:t>start,</w:t></w:r><w:r><w:rPr/><w:t xml:space="preserve"> especially with starting up a new repository in </w:t></w:r><w:r><w:rPr/><w:t>GitHub</w:t></w:r><w:r><w:rPr/><w:t xml:space="preserve"> to link the two together. But a</w:t></w:r><w:r><w:rPr/><w:t xml:space="preserve">s </w:t></w:r><w:r><w:rPr/><w:t>I</w:t></w:r><w:r><w:rPr/><w:t xml:space="preserve"> made it through the worksheets it started to become a lot easier</w:t></w:r><w:r><w:rPr/><w:t xml:space="preserve"> and </w:t></w:r><w:r><w:rPr/><w:t>I</w:t></w:r><w:r><w:rPr/><w:t xml:space="preserve"> started to understand </w:t></w:r><w:r><w:rPr/><w:t>how to create</w:t></w:r><w:r><w:rPr/><w:t xml:space="preserve"> a border, padding, </w:t></w:r><w:r><w:rPr/><w:t>flexbox</w:t></w:r><w:r><w:rPr/><w:t xml:space="preserve"> and other attributes to create a simple website</w:t></w:r><w:r><w:rPr/><w:t xml:space="preserve">. I did find it difficult when uploading images onto </w:t></w:r><w:r><w:rPr/><w:t>GitHub</w:t></w:r><w:r><w:rPr/><w:t xml:space="preserve"> as they </w:t></w:r><w:r><w:rPr/><w:t>didn''t</w:t></w:r><w:r><w:rPr/><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr/><w:t>always</w:t></w:r><w:r><w:rPr/><w:t xml:space="preserve"> work but </w:t></w:r><w:r><w:rPr/><w:t>I</w:t></w:r><w:r><w:rPr/><w:t xml:space="preserve"> managed to work</w:t></w:r><w:r><w:rPr/><w:t xml:space="preserve"> through this problem and finish the coding sheets anyway.</w:t></w:r></w:p><w:p/><w:p><w:r><w:rPr/><w:t xml:space="preserve">For the last couple </w:t></w:r><w:r><w:rPr/><w:t>worksheets</w:t></w:r><w:r><w:rPr/><w:t>,</w:t></w:r><w:r><w:rPr/><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr/><w:t xml:space="preserve">we looked at adding links to different websites, adding audio, </w:t></w:r><w:r><w:rPr/><w:t>anchor</w:t></w:r><w:r><w:rPr/><w:t xml:space="preserve"> links, nav bar, google and adobe font</w:t></w:r><w:r><w:rPr/><w:t>s.</w:t></w:r><w:r><w:rPr/><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr/><w:t xml:space="preserve">This was all new to me and </w:t></w:r><w:r><w:rPr/><w:t>i</w:t></w:r><w:r><w:rPr/><w:t xml:space="preserve"> found it interesting to find out how websites are </w:t></w:r><w:r><w:rPr/><w:t>made,</w:t></w:r><w:r><w:rPr/><w:t xml:space="preserve"> which </w:t></w:r><w:r><w:rPr/><w:t>I</w:t></w:r><w:r><w:rPr/><w:t xml:space="preserve"> found </w:t></w:r><w:r><w:rPr/><w:t>very useful</w:t></w:r><w:r><w:rPr/><w:t xml:space="preserve">. </w:t></w:r></w:p><w:p><w:r><w:rPr/><w:t>Overall,</w:t></w:r><w:r><w:rPr/><w:t xml:space="preserve"> the coding sheets were a good introduction to HTML and CSS </w:t></w:r><w:r><w:rPr/><w:t xml:space="preserve">and really helped when creating </w:t></w:r><w:r><w:rPr/><w:t>my</w:t></w:r><w:r><w:rPr/><w:t xml:space="preserve"> website. </w:t></w:r></w:p><w:p/>'
$r.InsertXML($xml)
